$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old merges
$ws.Range("C3:D3").UnMerge()
$ws.Range("C13:D13").UnMerge()

# Clear old content/format entirely within prior used range
$ws.Range("C3:D14").Clear()

# Row 4: long description header (no style)
$ws.Range("C4").Value2 = "Spreadsheet SpreadsheetResult spr(String name, Integer age)"

# Row 5: table header
$ws.Range("C5").Value2 = "Steps"
$ws.Range("D5").Value2 = "Formula"

# Row 6: Step1 + formula-as-text (quote prefix)
$ws.Range("C6").Value2 = "Step1"
$ws.Range("D6").Value = "'= sayHello(name)"

# Row 7: Step2 + formula-as-text (quote prefix)
$ws.Range("C7").Value2 = "Step2"
$ws.Range("D7").Value = "'= ""I am "" + age + "" age old."""

# Row 14: Environment
$ws.Range("C14").Value2 = "Environment"

# Row 15: dependency / Project2-*
$ws.Range("C15").Value2 = "dependency"
$ws.Range("D15").Value2 = "Project2-*"
